# Improve Excel exporter to create new sheets for new years
#
# This particular workbook snapshot reflects a re-export of the existing
# "2021" sheet: two of the previously recorded time entries (Implementierung
# and Maintenance, rows 6 and 7) were removed/reset, and two new entries
# for a new "Company" bucket (Administration/adm and Support/sup) were
# appended as rows 9 and 10. The running totals in row 3/4 are formulas and
# recalculate automatically once the underlying D-column values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Software / Implementierung / div) and row 7 (Software / Maintenance
# / test) lose their recorded minutes - the cells become empty again.
$ws.Cells.Item(6, 4).ClearContents()
$ws.Cells.Item(7, 4).ClearContents()

# New bucket "Company" with two sub-categories, appended below the existing
# "Software" rows. Write column C, then B, then A so the shared-string table
# is built up in the same order as the source edit (adm, Administration,
# Company, sup, Support), and so "Company" is reused (not duplicated) for
# the second row.
$ws.Cells.Item(9, 3).Value = "adm"
$ws.Cells.Item(9, 2).Value = "Administration"
$ws.Cells.Item(9, 1).Value = "Company"

$ws.Cells.Item(10, 3).Value = "sup"
$ws.Cells.Item(10, 2).Value = "Support"
$ws.Cells.Item(10, 1).Value = "Company"

# The active selection moved from D11 to D8.
$ws.Range("D8").Select()
